$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5586043333333334
$ws.Range("H2").Value = 1.675813
$ws.Range("I2").Value = 0.01643366487114074
$ws.Range("J2").Value = 0.01643366487114074
$ws.Range("M2").Value = 61.04160633333334
$ws.Range("N2").Value = 183.124819
$ws.Range("O2").Value = 0.2043613460574534
$ws.Range("P2").Value = 0.2043613460574534
$ws.Range("Q2").Value = 34.09810581142745
$ws.Range("R2").Value = 306.8829523028471
$ws.Range("S2").Value = 0.003358405873723409
$ws.Range("T2").Value = 0.003358405873723408
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.5586043333333334
$ws.Range("H3").Value = 1.675813
$ws.Range("I3").Value = 0.01643366487114074
$ws.Range("J3").Value = 0.01643366487114074
$ws.Range("O3").Value = 0.3559304658284363
$ws.Range("P3").Value = 0.3559304658284363
$ws.Range("Q3").Value = 59.38772140361934
$ws.Range("R3").Value = 534.489492632574
$ws.Range("S3").Value = 0.005849241992853533
$ws.Range("T3").Value = 0.005849241992853533
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5586043333333334
$ws.Range("H4").Value = 1.675813
$ws.Range("I4").Value = 0.01643366487114074
$ws.Range("J4").Value = 0.01643366487114074
$ws.Range("M4").Value = 131.3384093333333
$ws.Range("N4").Value = 394.015228
$ws.Range("O4").Value = 0.4397081881141102
$ws.Range("P4").Value = 0.4397081881141103
$ws.Range("Q4").Value = 73.36620458670711
$ws.Range("R4").Value = 660.295841280364
$ws.Range("S4").Value = 0.007226017004563798
$ws.Range("T4").Value = 0.007226017004563797
$ws.Range("H5").Value = 63.825936
$ws.Range("I5").Value = 0.6259016025719319
$ws.Range("J5").Value = 0.6259016025719319
$ws.Range("M5").Value = 61.04160633333334
$ws.Range("N5").Value = 183.124819
$ws.Range("O5").Value = 0.2043613460574534
$ws.Range("P5").Value = 0.2043613460574534
$ws.Range("Q5").Value = 1298.679219722843
$ws.Range("R5").Value = 11688.11297750559
$ws.Range("S5").Value = 0.1279100940011173
$ws.Range("T5").Value = 0.1279100940011173
$ws.Range("H6").Value = 63.825936
$ws.Range("I6").Value = 0.6259016025719319
$ws.Range("J6").Value = 0.6259016025719319
$ws.Range("O6").Value = 0.3559304658284363
$ws.Range("P6").Value = 0.3559304658284363
$ws.Range("S6").Value = 0.2227774489661925
$ws.Range("T6").Value = 0.2227774489661926
$ws.Range("H7").Value = 63.825936
$ws.Range("I7").Value = 0.6259016025719319
$ws.Range("J7").Value = 0.6259016025719319
$ws.Range("M7").Value = 131.3384093333333
$ws.Range("N7").Value = 394.015228
$ws.Range("O7").Value = 0.4397081881141102
$ws.Range("P7").Value = 0.4397081881141103
$ws.Range("Q7").Value = 2794.265636150378
$ws.Range("R7").Value = 25148.39072535341
$ws.Range("S7").Value = 0.2752140596046221
$ws.Range("T7").Value = 0.2752140596046221
$ws.Range("G8").Value = 12.157548
$ws.Range("H8").Value = 36.472644
$ws.Range("I8").Value = 0.3576647325569273
$ws.Range("J8").Value = 0.3576647325569273
$ws.Range("M8").Value = 61.04160633333334
$ws.Range("N8").Value = 183.124819
$ws.Range("O8").Value = 0.2043613460574534
$ws.Range("P8").Value = 0.2043613460574534
$ws.Range("Q8").Value = 742.1162589946041
$ws.Range("R8").Value = 6679.046330951438
$ws.Range("S8").Value = 0.07309284618261276
$ws.Range("T8").Value = 0.07309284618261276
$ws.Range("G9").Value = 12.157548
$ws.Range("H9").Value = 36.472644
$ws.Range("I9").Value = 0.3576647325569273
$ws.Range("J9").Value = 0.3576647325569273
$ws.Range("O9").Value = 0.3559304658284363
$ws.Range("P9").Value = 0.3559304658284363
$ws.Range("Q9").Value = 1292.523223489368
$ws.Range("R9").Value = 11632.70901140431
$ws.Range("S9").Value = 0.1273037748693902
$ws.Range("T9").Value = 0.1273037748693902
$ws.Range("G10").Value = 12.157548
$ws.Range("H10").Value = 36.472644
$ws.Range("I10").Value = 0.3576647325569273
$ws.Range("J10").Value = 0.3576647325569273
$ws.Range("M10").Value = 131.3384093333333
$ws.Range("N10").Value = 394.015228
$ws.Range("O10").Value = 0.4397081881141102
$ws.Range("P10").Value = 0.4397081881141103
$ws.Range("Q10").Value = 1596.753015713648
$ws.Range("R10").Value = 14370.77714142283
$ws.Range("S10").Value = 0.1572681115049243
$ws.Range("T10").Value = 0.1572681115049243
